# Updated cryptos list with latest prices / 1h volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.401.38'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').Value = '1.937.15'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('D4').Value = '''1.000'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = '''0.7709'
$ws.Range('E5').Value = '  +8.30%  '
$ws.Range('D6').Value = '''248.27'
$ws.Range('E6').Value = '  -1.37%  '
$ws.Range('D7').Value = '''1.000'
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('D8').Value = '''27.95'
$ws.Range('E8').Value = '  +2.34%  '
$ws.Range('D9').Value = '''0.3199'
$ws.Range('E9').Value = '  -3.33%  '
$ws.Range('D10').Value = '''0.07102'
$ws.Range('E10').Value = '  -3.28%  '
$ws.Range('D11').Value = '''0.7837'
$ws.Range('E11').Value = '  -2.83%  '
$ws.Range('D12').Value = '''0.08012'
$ws.Range('E12').Value = '  -0.84%  '
$ws.Range('D13').Value = '1.939.15'
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').Value = '''5.377'
$ws.Range('E14').Value = '  -2.27%  '
$ws.Range('D15').Value = '''95.11'
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('D16').Value = '''14.54'
$ws.Range('E16').Value = '  -4.03%  '
$ws.Range('D17').Value = '30.399.65'
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D18').Value = '''257.64'
$ws.Range('E18').Value = '  +1.82%  '
$ws.Range('D19').Value = '''0.000008017'
$ws.Range('E19').Value = '  -2.54%  '
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('D21').Value = '2.191.95'
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').Value = '''1.000'
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').Value = '''1.000'
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('D24').Value = '''6.762'
$ws.Range('E24').Value = '  -3.62%  '
$ws.Range('D25').Value = '''9.621'
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('D26').Value = '''164.53'
$ws.Range('E26').Value = '  +0.58%  '
$ws.Range('D27').Value = '''19.16'
$ws.Range('E27').Value = '  -0.95%  '
$ws.Range('D28').Value = '''0.1339'
$ws.Range('E28').Value = '  +2.65%  '
$ws.Range('D29').Value = '''2.294'
$ws.Range('E29').Value = '  -2.64%  '
$ws.Range('E30').Value = '  +1.15%  '
$ws.Range('D31').Value = '''1.530'
$ws.Range('E31').Value = '  -2.77%  '
$ws.Range('D32').Value = '''4.439'
$ws.Range('E32').Value = '  +0.36%  '
$ws.Range('D33').Value = '''4.159'
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('D34').Value = '''0.05196'
$ws.Range('E34').Value = '  +0.11%  '
$ws.Range('D35').Value = '''1.282'
$ws.Range('E35').Value = '  +0.96%  '
$ws.Range('D36').Value = '''0.7508'
$ws.Range('E36').Value = '  +0.50%  '
$ws.Range('D37').Value = '''2.778'
$ws.Range('E37').Value = '  +0.93%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('E39').Value = '  -0.27%  '
$ws.Range('D40').Value = '''78.37'
$ws.Range('E40').Value = '  -0.81%  '
$ws.Range('D41').Value = '''6.465'
$ws.Range('E41').Value = '  +0.74%  '
$ws.Range('D42').Value = '''0.4525'
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('D43').Value = '''1.978'
$ws.Range('E43').Value = '  -1.75%  '
$ws.Range('D44').Value = '''1.001'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').Value = '''0.8357'
$ws.Range('E45').Value = '  -1.52%  '
$ws.Range('D46').Value = '''101.47'
$ws.Range('E46').Value = '  -0.21%  '
$ws.Range('D47').Value = '''9.859'
$ws.Range('E47').Value = '  +1.92%  '
$ws.Range('D48').Value = '''7.519'
$ws.Range('E48').Value = '  +0.69%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '''37.53'
$ws.Range('E49').Value = '  +2.29%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '''985.70'
$ws.Range('E50').Value = '  +11.15%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '''1.503'
$ws.Range('E51').Value = '  +1.17%  '
